# Update a few data values in the training schedule sheet and leave the
# selection on the cell that was last edited (E2), matching the authored
# workbook's saved UI state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 7
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 13

$ws.Range("E2").Select() | Out-Null
